# Delete the two "spike in" library rows for sample 26a/26b (rows 27 and 28),
# shifting the remaining rows up, as described in the commit message:
# "I deleted 2 of the libraries in fastq03.10.20 ... (26a/b) ... to simplify the database."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 27 and 28 entirely (whole-row delete shifts everything below up by 2).
$ws.Rows("27:28").Delete()

# Update the selection to match the saved state after the edit.
$ws.Range("C26").Select()
